$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suite 1")

# Row 3 describes a test step: Object = FileSystem, Action = checkFileContainsKeyword.
# The old "gaga" placeholder input is no longer needed, so both the
# Description and Input cells are cleared.
$ws.Range("B3").Value = "FileSystem"
$ws.Range("C3").Value = "checkFileContainsKeyword"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
